$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 27 (the old blank "F27/G27" row),
# which shifts rows 27-34 down to 28-35 and keeps their formatting intact.
$ws.Rows("27:27").Insert()

# Populate the newly inserted row 27 with the new field mapping:
# NotificationMainLotMainProductCode
$ws.Range("A27").Value = "notifications"
$ws.Range("B27").Value = "notifications_<region>"
$ws.Range("C27").Value = "/*/*/oos:lots/oos:lot/oos:products/oos:product/oos:code"
$ws.Range("D27").Value = "oos:lots/oos:lot/oos:products/oos:product/oos:code"
$ws.Range("E27").Value = '"2320000"'
$ws.Range("F27").Value = "factor"
$ws.Range("G27").Value = "NotificationMainLotMainProductCode"

# Restore the view state as closely as possible: scroll position (top-left
# visible cell under the frozen pane) and active cell / selection.
$excel.ActiveWindow.ScrollRow = 14
$ws.Range("D21").Select()
